$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1660.0625
$ws.Range("I32").Value = 1833.3334
$ws.Range("J32").Value = 1620.0769
$ws.Range("K32").Value = 1833.3334
$ws.Range("L32").Value = 1620.0769
$ws.Range("M32").Value = -1507.3334
$ws.Range("N32").Value = -2272.0769
$ws.Range("H64").Value = 4223.25
$ws.Range("I64").Value = 4492
$ws.Range("J64").Value = 4133.6665
$ws.Range("K64").Value = 4492
$ws.Range("L64").Value = 4133.6665
$ws.Range("M64").Value = -4244
$ws.Range("N64").Value = -4629.6665
$ws.Range("H67").Value = 4223.25
$ws.Range("I67").Value = 4492
$ws.Range("J67").Value = 4133.6665
$ws.Range("K67").Value = 4492
$ws.Range("L67").Value = 4133.6665
$ws.Range("M67").Value = -3634
$ws.Range("N67").Value = -5849.6665
$ws.Range("H74").Value = 4926.6665
$ws.Range("I74").Value = 4760
$ws.Range("J74").Value = 5135
$ws.Range("K74").Value = 4760
$ws.Range("L74").Value = 5135
$ws.Range("M74").Value = -3824
$ws.Range("N74").Value = -7007
$ws.Range("H76").Value = 3675.6667
$ws.Range("I76").Value = 3301
$ws.Range("J76").Value = 3722.5
$ws.Range("K76").Value = 3301
$ws.Range("L76").Value = 3722.5
$ws.Range("M76").Value = -2986
$ws.Range("N76").Value = -4352.5
$ws.Range("H77").Value = 4926.6665
$ws.Range("I77").Value = 4760
$ws.Range("J77").Value = 5135
$ws.Range("K77").Value = 23800
$ws.Range("L77").Value = 25675
$ws.Range("M77").Value = -19120
$ws.Range("N77").Value = -35035
$ws.Range("H79").Value = 3675.6667
$ws.Range("I79").Value = 3301
$ws.Range("J79").Value = 3722.5
$ws.Range("K79").Value = 3301
$ws.Range("L79").Value = 3722.5
$ws.Range("M79").Value = -2209
$ws.Range("N79").Value = -5906.5
$ws.Range("H95").Value = 32824.75
$ws.Range("J95").Value = 32824.75
$ws.Range("L95").Value = 32824.75
$ws.Range("N95").Value = -38316.75
$ws.Range("H125").Value = 3100
$ws.Range("I125").Value = 2466.6667
$ws.Range("J125").Value = 5000
$ws.Range("K125").Value = 22200.0003
$ws.Range("L125").Value = 45000
$ws.Range("M125").Value = -19740.0003
$ws.Range("N125").Value = -49920
$ws.Range("H137").Value = 1826.6945
$ws.Range("I137").Value = 2189.0476
$ws.Range("J137").Value = 1319.4
$ws.Range("K137").Value = 6567.1428
$ws.Range("L137").Value = 3958.2
$ws.Range("M137").Value = -4017.1428
$ws.Range("N137").Value = -9058.200000000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2947.0688
$ws.Range("I63").Value = 2100.2778
$ws.Range("J63").Value = 4332.727
$ws.Range("K63").Value = 2100.2778
$ws.Range("L63").Value = 4332.727
$ws.Range("M63").Value = -1414.2778
$ws.Range("N63").Value = -5704.727
$ws.Range("H66").Value = 2947.0688
$ws.Range("I66").Value = 2100.2778
$ws.Range("J66").Value = 4332.727
$ws.Range("K66").Value = 10501.389
$ws.Range("L66").Value = 21663.635
$ws.Range("M66").Value = -7069.388999999999
$ws.Range("N66").Value = -28527.635
$ws.Range("H88").Value = 1925
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 1900
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 1900
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -2712
$ws.Range("H91").Value = 1925
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 1900
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 1900
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -4708
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1230.8889
$ws.Range("J86").Value = 1162.5
$ws.Range("L86").Value = 1162.5
$ws.Range("N86").Value = -3408.5
$ws.Range("H89").Value = 1230.8889
$ws.Range("J89").Value = 1162.5
$ws.Range("L89").Value = 5812.5
$ws.Range("N89").Value = -17044.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5292915
$ws.Range("I31").Value = 1431
$ws.Range("J31").Value = 22225664
$ws.Range("K31").Value = 1431
$ws.Range("L31").Value = 22225664
$ws.Range("M31").Value = -1136
$ws.Range("N31").Value = -22226254
$ws.Range("H34").Value = 5292915
$ws.Range("I34").Value = 1431
$ws.Range("J34").Value = 22225664
$ws.Range("K34").Value = 1431
$ws.Range("L34").Value = 22225664
$ws.Range("M34").Value = -1229
$ws.Range("N34").Value = -22226068
$ws.Range("H62").Value = 9210.833000000001
$ws.Range("I62").Value = 2590.4546
$ws.Range("J62").Value = 19614.285
$ws.Range("K62").Value = 2590.4546
$ws.Range("L62").Value = 19614.285
$ws.Range("M62").Value = -1966.4546
$ws.Range("N62").Value = -20862.285
$ws.Range("H65").Value = 9210.833000000001
$ws.Range("I65").Value = 2590.4546
$ws.Range("J65").Value = 19614.285
$ws.Range("K65").Value = 12952.273
$ws.Range("L65").Value = 98071.425
$ws.Range("M65").Value = -9832.273000000001
$ws.Range("N65").Value = -104311.425
$ws.Range("H134").Value = 2017.4517
$ws.Range("I134").Value = 870.3182
$ws.Range("J134").Value = 4821.5557
$ws.Range("K134").Value = 2610.9546
$ws.Range("L134").Value = 14464.6671
$ws.Range("M134").Value = -75.95460000000003
$ws.Range("N134").Value = -19534.6671
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2852.861
$ws.Range("I132").Value = 988.4231
$ws.Range("J132").Value = 7700.4
$ws.Range("K132").Value = 8895.8079
$ws.Range("L132").Value = 69303.59999999999
$ws.Range("M132").Value = -6365.8079
$ws.Range("N132").Value = -74363.59999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 3850500
$ws.Range("I24").Value = 4620000
$ws.Range("J24").Value = 3000
$ws.Range("K24").Value = 4620000
$ws.Range("L24").Value = 3000
$ws.Range("M24").Value = -4619827
$ws.Range("N24").Value = -3346
$ws.Range("H70").Value = 6661.8
$ws.Range("I70").Value = 7166.6665
$ws.Range("J70").Value = 5904.5
$ws.Range("K70").Value = 7166.6665
$ws.Range("L70").Value = 5904.5
$ws.Range("M70").Value = -6896.6665
$ws.Range("N70").Value = -6444.5
$ws.Range("H73").Value = 6661.8
$ws.Range("I73").Value = 7166.6665
$ws.Range("J73").Value = 5904.5
$ws.Range("K73").Value = 7166.6665
$ws.Range("L73").Value = 5904.5
$ws.Range("M73").Value = -6230.6665
$ws.Range("N73").Value = -7776.5
$ws.Range("H80").Value = 2807.9473
$ws.Range("I80").Value = 2754.0625
$ws.Range("J80").Value = 3095.3333
$ws.Range("K80").Value = 2754.0625
$ws.Range("L80").Value = 3095.3333
$ws.Range("M80").Value = -1756.0625
$ws.Range("N80").Value = -5091.3333
$ws.Range("H83").Value = 2807.9473
$ws.Range("I83").Value = 2754.0625
$ws.Range("J83").Value = 3095.3333
$ws.Range("K83").Value = 13770.3125
$ws.Range("L83").Value = 15476.6665
$ws.Range("M83").Value = -8778.3125
$ws.Range("N83").Value = -25460.6665
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = ""
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = ""
$ws.Range("H132").Value = 43481028
$ws.Range("I132").Value = 100002710
$ws.Range("J132").Value = 2807.5386
$ws.Range("K132").Value = 300008130
$ws.Range("L132").Value = 8422.6158
$ws.Range("M132").Value = -300005600
$ws.Range("N132").Value = -13482.6158
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2454.3635
$ws.Range("I132").Value = 1887.25
$ws.Range("J132").Value = 2988.1177
$ws.Range("K132").Value = 5661.75
$ws.Range("L132").Value = 8964.3531
$ws.Range("M132").Value = -3131.75
$ws.Range("N132").Value = -14024.3531
$ws.Range("H136").Value = 6946358
$ws.Range("I136").Value = 14707293
$ws.Range("J136").Value = 2362.9473
$ws.Range("K136").Value = 44121879
$ws.Range("L136").Value = 7088.841899999999
$ws.Range("M136").Value = -44119329
$ws.Range("N136").Value = -12188.8419
